# Ardi y Stephen: Se arreglo presentación
#
# 1) Bump the cached "datetimeFigureOut" field text from 01/10/2019 to
#    02/10/2019 everywhere it is cached: the Notes Master, the Slide
#    Master, and slide layouts 7-13 (the only layouts/masters that carry
#    that cached field text).
# 2) Slide 2: widen + rename the "CuadroTexto 14" caption from
#    "Diccionario de datos" to "Diagrama de despliegue".

$p = $ppt.ActivePresentation

$newDate = "02/10/2019"

# --- Slide Master: "Marcador de fecha 3" date placeholder ---
$design = $p.Designs.Item(1)
$slideMaster = $design.SlideMaster
$slideMaster.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

# --- Slide Layouts 7-13: each layout's own date placeholder shape ---
$layoutDateShapeIndex = @{
    7  = 6
    8  = 2
    9  = 1
    10 = 4
    11 = 4
    12 = 3
    13 = 3
}

foreach ($layoutIdx in $layoutDateShapeIndex.Keys) {
    $layout = $slideMaster.CustomLayouts.Item($layoutIdx)
    $shapeIdx = $layoutDateShapeIndex[$layoutIdx]
    $layout.Shapes.Item($shapeIdx).TextFrame.TextRange.Text = $newDate
}

# --- Slide 2: resize + rename "CuadroTexto 14" ---
$s2 = $p.Slides.Item(2)
$caption = $s2.Shapes.Item(11)
$caption.Width = 211.4449606
$caption.TextFrame.TextRange.Text = "Diagrama de despliegue"
